# Updates cryptos list values (Price / Volume(1h)) to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.003.72'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '1.742.32'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('D4').Value = "'0.9995"
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = "'247.42"
$ws.Range('E5').Value = '  +3.18%  '
$ws.Range('D6').Value = "'0.9999"
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').Value = "'0.5048"
$ws.Range('E7').Value = '  -4.40%  '
$ws.Range('D8').Value = "'0.2745"
$ws.Range('E8').Value = '  +0.44%  '
$ws.Range('D9').Value = "'0.06181"
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('D10').Value = '1.751.47'
$ws.Range('D11').Value = "'0.07258"
$ws.Range('E11').Value = '  +1.04%  '
$ws.Range('D12').Value = "'0.6546"
$ws.Range('E12').Value = '  +2.08%  '
$ws.Range('D13').Value = "'15.13"
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('D14').Value = "'4.681"
$ws.Range('E14').Value = '  +1.69%  '
$ws.Range('D15').Value = "'77.64"
$ws.Range('E15').Value = '  +0.17%  '
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').Value = "'0.9994"
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').Value = '26.017.93'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').Value = "'11.92"
$ws.Range('E19').Value = '  +1.27%  '
$ws.Range('D20').Value = "'0.000006856"
$ws.Range('E20').Value = '  +1.64%  '
$ws.Range('D21').Value = '1.977.09'
$ws.Range('E21').Value = '  +0.66%  '
$ws.Range('D22').Value = "'4.470"
$ws.Range('E22').Value = '  +2.63%  '
$ws.Range('D23').Value = "'8.720"
$ws.Range('E23').Value = '  +1.19%  '
$ws.Range('D24').Value = "'5.382"
$ws.Range('E24').Value = '  +2.57%  '
$ws.Range('D25').Value = "'135.56"
$ws.Range('E25').Value = '  -3.26%  '
$ws.Range('D26').Value = "'1.504"
$ws.Range('E26').Value = '  -0.73%  '
$ws.Range('D27').Value = "'15.26"
$ws.Range('E27').Value = '  +0.27%  '
$ws.Range('D28').Value = "'1.784"
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('D29').Value = "'105.38"
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('D30').Value = "'3.924"
$ws.Range('E30').Value = '  +2.33%  '
$ws.Range('D31').Value = "'0.08167"
$ws.Range('E31').Value = '  -2.69%  '
$ws.Range('D32').Value = "'3.681"
$ws.Range('E32').Value = '  +0.93%  '
$ws.Range('D33').Value = "'0.04683"
$ws.Range('E33').Value = '  +2.18%  '
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('D35').Value = "'0.9976"
$ws.Range('E35').Value = '  +0.62%  '
$ws.Range('E36').Value = '  -1.77%  '
$ws.Range('D37').Value = "'2.759"
$ws.Range('E37').Value = '  +2.08%  '
$ws.Range('D38').Value = "'0.01624"
$ws.Range('E38').Value = '  +1.33%  '
$ws.Range('D39').Value = "'1.930"
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').Value = "'100.91"
$ws.Range('E41').Value = '  +2.14%  '
$ws.Range('D42').Value = "'0.3921"
$ws.Range('E42').Value = '  +0.95%  '
$ws.Range('D43').Value = "'0.7630"
$ws.Range('E43').Value = '  +1.82%  '
$ws.Range('D44').Value = "'5.015"
$ws.Range('E44').Value = '  +1.39%  '
$ws.Range('E45').Value = '  +1.45%  '
$ws.Range('D46').Value = "'6.323"
$ws.Range('E46').Value = '  +1.78%  '
$ws.Range('D47').Value = "'55.54"
$ws.Range('E47').Value = '  +1.51%  '
$ws.Range('D48').Value = "'0.05299"
$ws.Range('E48').Value = '  -0.22%  '
$ws.Range('D49').Value = "'30.71"
$ws.Range('E49').Value = '  -0.23%  '
$ws.Range('D50').Value = "'0.3474"
$ws.Range('E50').Value = '  +0.92%  '
$ws.Range('D51').Value = "'7.589"
$ws.Range('E51').Value = '  +0.89%  '
